# Updated symbol list on Fri Feb 17 06:32:07 UTC 2023 with GitHub Actions
# Applies refreshed price / 1h-volume (and, for rows 11-12, a coin-rank swap)
# quotes from the upstream coinranking.com scrape to the "cryptos" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: refreshed price/volume
$ws.Range("D2").Value = "'308.07"
$ws.Range("E2").Value = "'-4.56%"

# Row 3: refreshed price/volume
$ws.Range("D3").Value = "'49.40"
$ws.Range("E3").Value = "'-1.91%"

# Row 4: refreshed price/volume
$ws.Range("D4").Value = "'5.167"
$ws.Range("E4").Value = "'-3.33%"

# Row 5: refreshed price/volume
$ws.Range("D5").Value = "'0.07730"
$ws.Range("E5").Value = "'-5.34%"

# Row 6: refreshed price/volume
$ws.Range("E6").Value = "'-1.88%"

# Row 7: refreshed price/volume
$ws.Range("D7").Value = "'1.367"
$ws.Range("E7").Value = "'13.57%"

# Row 8: refreshed price/volume
$ws.Range("E8").Value = "'-7.25%"

# Row 9: refreshed price/volume
$ws.Range("D9").Value = "'0.1227"
$ws.Range("E9").Value = "'-8.47%"

# Row 10: refreshed price/volume
$ws.Range("D10").Value = "'0.1937"
$ws.Range("E10").Value = "'-1.23%"

# Row 11: coin-rank swap + refreshed price/volume
$ws.Range("B11").Value = "'MandalaExchangeToken"
$ws.Range("C11").Value = "'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D11").Value = "'0.09425"
$ws.Range("E11").Value = "'-2.52%"

# Row 12: coin-rank swap + refreshed price/volume
$ws.Range("B12").Value = "'BitrueCoin"
$ws.Range("C12").Value = "'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D12").Value = "'0.04642"
$ws.Range("E12").Value = "'3.79%"

# Row 13: refreshed price/volume
$ws.Range("D13").Value = "'0.1045"
$ws.Range("E13").Value = "'-0.29%"

# Row 14: refreshed price/volume
$ws.Range("D14").Value = "'0.001266"
$ws.Range("E14").Value = "'-3.96%"

# Row 15: refreshed price/volume
$ws.Range("D15").Value = "'0.04185"
$ws.Range("E15").Value = "'-2.59%"

# Row 16: refreshed price/volume
$ws.Range("D16").Value = "'0.005848"
$ws.Range("E16").Value = "'-0.73%"

# Row 17: refreshed price/volume
$ws.Range("E17").Value = "'2,023.59%"

# Row 18: refreshed price/volume
$ws.Range("D18").Value = "'3.335"
$ws.Range("E18").Value = "'-1.58%"

# Row 19: refreshed price/volume
$ws.Range("D19").Value = "'2.237"
$ws.Range("E19").Value = "'-8.26%"

# Row 20: refreshed price/volume
$ws.Range("D20").Value = "'0.3487"
$ws.Range("E20").Value = "'2.74%"

# Row 21: refreshed price/volume
$ws.Range("D21").Value = "'7.926"
$ws.Range("E21").Value = "'-2.71%"

# Row 22: refreshed price/volume
$ws.Range("D22").Value = "'0.1340"
$ws.Range("E22").Value = "'-5.50%"

# Row 23: refreshed price/volume
$ws.Range("D23").Value = "'0.3040"
$ws.Range("E23").Value = "'-0.30%"

# Row 24: refreshed price/volume
$ws.Range("D24").Value = "'0.001272"
$ws.Range("E24").Value = "'-2.41%"

# Row 25: refreshed price/volume
$ws.Range("D25").Value = "'0.003996"
$ws.Range("E25").Value = "'-6.41%"

# Row 26: refreshed price/volume
$ws.Range("D26").Value = "'0.0001352"
$ws.Range("E26").Value = "'0.27%"

# Row 38: refreshed price/volume
$ws.Range("D38").Value = "'0.02571"
$ws.Range("E38").Value = "'-7.02%"

# Row 39: refreshed price/volume
$ws.Range("D39").Value = "'0.05813"
$ws.Range("E39").Value = "'3.87%"

# Row 40: refreshed price/volume
$ws.Range("D40").Value = "'0.01074"
$ws.Range("E40").Value = "'70.70%"

# Row 41: refreshed price/volume
$ws.Range("D41").Value = "'0.007917"
$ws.Range("E41").Value = "'2.61%"

# Row 42: refreshed price/volume
$ws.Range("D42").Value = "'0.1421"
$ws.Range("E42").Value = "'-1.95%"

# Row 43: refreshed price/volume
$ws.Range("D43").Value = "'0.008441"
$ws.Range("E43").Value = "'9.98%"

# Row 44: refreshed price/volume
$ws.Range("D44").Value = "'0.007675"
$ws.Range("E44").Value = "'-4.96%"

# Row 45: refreshed price/volume
$ws.Range("D45").Value = "'0.3374"
$ws.Range("E45").Value = "'-3.80%"

# Row 46: refreshed price/volume
$ws.Range("D46").Value = "'0.00007021"
$ws.Range("E46").Value = "'3.00%"

# Row 47: refreshed price/volume
$ws.Range("E47").Value = "'0.34%"

# Row 48: refreshed price/volume
$ws.Range("D48").Value = "'0.05314"
$ws.Range("E48").Value = "'-13.34%"

# Row 49: refreshed price/volume
$ws.Range("E49").Value = "'0.23%"

# Row 50: refreshed price/volume
$ws.Range("E50").Value = "'0.34%"

# Row 51: refreshed price/volume
$ws.Range("E51").Value = "'0.34%"
